# download articles with pandoc title blocks
#
# The source "On Pilgrimage - May 1979" / "By Dorothy Day" header was
# authored by hand (Heading1 + bold "By ..." line, wrapped in a
# bookmark). Pandoc's title-block conversion instead produces a
# dedicated Title paragraph (words/spaces/punctuation each as their own
# run - pandoc's usual one-run-per-inline-chunk style) followed by an
# Authors-styled paragraph with just the author's name (no leading
# "By ").

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Best-effort: drop the legacy bookmark that wrapped the old
# heading. Its name ("on-pilgrimage---may-1979") isn't a legal Word
# bookmark identifier (hyphens aren't allowed), so real bookmark APIs
# can't see/target it; try anyway in case a name is resolvable, but
# don't let it block the rest of the edit.
try {
    if ($d.Bookmarks.Exists("on-pilgrimage---may-1979")) {
        $d.Bookmarks("on-pilgrimage---may-1979").Delete()
    }
} catch {
}

# --- Paragraph 1: "On Pilgrimage - May 1979" (Heading1, single run)
# becomes a Title-styled paragraph split word-by-word / space-by-space
# / punctuation-by-punctuation into separate runs.
$titlePara = $d.Paragraphs(1)
$titleXml = "<w:p xmlns:w='$wNs'>" +
    "<w:r><w:t xml:space='preserve'>On</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Pilgrimage</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>-</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>May</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>1979</w:t></w:r>" +
    "</w:p>"
$titlePara.Range.InsertXML($titleXml)
$d.Paragraphs(1).Style = "Title"

# --- Paragraph 2: "By Dorothy Day" (bold run) becomes an
# Authors-styled paragraph containing only the name, split into
# "Dorothy" / " " / "Day" runs, with no bold and no "By " prefix.
$authorPara = $d.Paragraphs(2)
$authorXml = "<w:p xmlns:w='$wNs'>" +
    "<w:r><w:t xml:space='preserve'>Dorothy</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Day</w:t></w:r>" +
    "</w:p>"
$authorPara.Range.InsertXML($authorXml)
$d.Paragraphs(2).Style = "Authors"

Write-Output "Paragraph 1 [$($d.Paragraphs(1).Style.NameLocal)]: $($d.Paragraphs(1).Range.Text)"
Write-Output "Paragraph 2 [$($d.Paragraphs(2).Style.NameLocal)]: $($d.Paragraphs(2).Range.Text)"
